# authorization_reports -> distribution reports. translations individual -> member.
# More table column name changes.
#
# This script reproduces, via Excel COM automation, an editing session in
# which the author:
#   1. Was on the "model" sheet and renamed the `individual_id` field to
#      `member_id` (row 7, column B), then widened columns A and B so the
#      longer field names are fully visible.
#   2. Switched to the "table_specific_translations" sheet and renamed the
#      matching translation label `Individual ID` -> `Member ID` (row 7,
#      column B), leaving that sheet focused/selected at B7.

$wb = $excel.ActiveWorkbook

# --- 1. "model" sheet: individual_id -> member_id -------------------------
$model = $wb.Worksheets.Item("model")
$model.Activate()
$model.Range("B7").Value = "member_id"

# Widen columns A and B so the (now longer) field names fit, matching the
# <cols> block added to this sheet.
$model.Columns.Item(1).ColumnWidth = 42.833333333333336
$model.Columns.Item(2).ColumnWidth = 40

# Restore the selection this sheet had before focus moved away.
$model.Range("B12").Select()

# --- 2. "table_specific_translations" sheet: individual_id -> member_id ---
$translations = $wb.Worksheets.Item("table_specific_translations")
$translations.Activate()
$translations.Range("A7").Value = "member_id"
$translations.Range("B7").Value = "Member ID"
$translations.Range("B7").Select()
